# Update the COVID-19 Valais figures workbook with newly reported data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L ("Nb nouveaux deces a l'hopital") and M ("Nb nouveaux deces
# extra-hospitaliers") carry a Text ("@") number format even though the
# sheet always stores genuine numeric values in them. Writing straight to
# .Value on a Text-formatted cell stores the input as a string, so round
# trip the format to General, assign the number, then restore "@" so the
# cell keeps its original numeric type and style.
function Set-NumericValue($addr, $val) {
    $rng = $ws.Range($addr)
    $savedFormat = $rng.NumberFormat()
    $rng.NumberFormat = "General"
    $rng.Value = $val
    $rng.NumberFormat = $savedFormat
}

# --- Row 359 (A359=44245): revised new-case count
$ws.Range("C359").Value = 45

# --- Row 360 (A360=44246): revised new-case count and new extra-hospital death
$ws.Range("C360").Value = 49
Set-NumericValue "M360" 2

# --- Row 361 (A361=44247): new extra-hospital death recorded
Set-NumericValue "M361" 1

# --- Row 362 (A362=44248): revised new-case count and new extra-hospital death
$ws.Range("C362").Value = 36
Set-NumericValue "M362" 1

# --- Row 363 (A363=44249): revised new-case count
$ws.Range("C363").Value = 69

# --- Row 364 (A364=44250): brand-new day of data added
$ws.Range("C364").Value = 17
$ws.Range("E364").Value = 10
$ws.Range("F364").Value = 8
$ws.Range("G364").Value = 32
Set-NumericValue "L364" 0
Set-NumericValue "M364" 1

# --- Restore the active selection on the frozen (bottom-right) pane to A2
$ws.Range("A2").Select()
